$wb = $excel.ActiveWorkbook

# --- AddTest sheet: clear the stale single-cell selection, select the used range instead ---
$addTest = $wb.Worksheets.Item("AddTest")
[void]$addTest.Range("A1:C2").Select()

# --- Insert a new sheet "tc048" right after "AddTest", extending the AddTest data with
#     extra tcname/tcdescription/priority/QA columns ---
$newSheet = $wb.Worksheets.Add($null, $addTest)
$newSheet.Name = "tc048"

$newSheet.Range("A1").Value = "Epic"
$newSheet.Range("A2").Value = "Epic Mohit"
$newSheet.Range("B1").Value = "Feature"
$newSheet.Range("B2").Value = "Mohit Feature"
$newSheet.Range("C1").Value = "Requirement"
$newSheet.Range("C2").Value = "RQ-489"
$newSheet.Range("D1").Value = "Tcname"
$newSheet.Range("D2").Value = "Unit testing "
$newSheet.Range("E1").Value = "Tcdescription"
$newSheet.Range("E2").Value = "work"
$newSheet.Range("F1").Value = "priority"
$newSheet.Range("F2").Value = "Low"
$newSheet.Range("G1").Value = "QA"
$newSheet.Range("G2").Value = "Mohit Aman"

$newSheet.Range("A1:G2").WrapText = $true
[void]$newSheet.Range("G8").Select()
